$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 274.75
$ws.Range("I5").Value = 166.33333
$ws.Range("J5").Value = 600
$ws.Range("K5").Value = 166.33333
$ws.Range("L5").Value = 600
$ws.Range("M5").Value = -51.33332999999999
$ws.Range("N5").Value = -830

$ws.Range("H8").Value = 200
$ws.Range("I8").Value = 200
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 600
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -461
$ws.Range("N8").ClearContents()

$ws.Range("H40").Value = 2160
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H75").Value = 15313.5
$ws.Range("J75").Value = 15313.5
$ws.Range("L75").Value = 15313.5
$ws.Range("N75").Value = -17185.5

$ws.Range("H78").Value = 15313.5
$ws.Range("J78").Value = 15313.5
$ws.Range("L78").Value = 45940.5
$ws.Range("N78").Value = -55300.5

$ws.Range("H80").Value = 355.42856
$ws.Range("I80").Value = 195.45454
$ws.Range("J80").Value = 531.4
$ws.Range("K80").Value = 586.3636200000001
$ws.Range("L80").Value = 1594.2
$ws.Range("M80").Value = 411.6363799999999
$ws.Range("N80").Value = -3590.2

$ws.Range("H83").Value = 355.42856
$ws.Range("I83").Value = 195.45454
$ws.Range("J83").Value = 531.4
$ws.Range("K83").Value = 1759.09086
$ws.Range("L83").Value = 4782.599999999999
$ws.Range("M83").Value = 3232.90914
$ws.Range("N83").Value = -14766.6

$ws.Range("H112").Value = 4400
$ws.Range("J112").Value = 4400
$ws.Range("L112").Value = 13200
$ws.Range("N112").Value = -15416

$ws.Range("H127").Value = 662.4167
$ws.Range("I127").Value = 353.25
$ws.Range("K127").Value = 1059.75
$ws.Range("M127").Value = 3900.25

$ws.Range("H128").Value = 34899.5
$ws.Range("J128").Value = 34899.5
$ws.Range("L128").Value = 34899.5
$ws.Range("N128").Value = -44859.5

$ws.Range("H132").Value = 10105597
$ws.Range("I132").Value = 10755177
$ws.Range("J132").Value = 37100
$ws.Range("K132").Value = 32265531
$ws.Range("L132").Value = 111300
$ws.Range("M132").Value = -32263001
$ws.Range("N132").Value = -116360

$ws.Range("H137").Value = 1447.4
$ws.Range("I137").Value = 1375.4166
$ws.Range("J137").Value = 1735.3334
$ws.Range("K137").Value = 4126.2498
$ws.Range("L137").Value = 5206.0002
$ws.Range("M137").Value = -1576.2498
$ws.Range("N137").Value = -10306.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 990
$ws.Range("J4").Value = 990
$ws.Range("L4").Value = 990
$ws.Range("N4").Value = -1222

$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

$ws.Range("H32").Value = 7826
$ws.Range("I32").Value = 7523.6665
$ws.Range("K32").Value = 7523.6665
$ws.Range("M32").Value = -7236.6665

$ws.Range("H45").Value = 2584.5
$ws.Range("I45").Value = 2584.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2584.5
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2207.5
$ws.Range("N45").ClearContents()

$ws.Range("H61").Value = 1449.8889
$ws.Range("I61").Value = 1141.5
$ws.Range("J61").Value = 2066.6667
$ws.Range("K61").Value = 1141.5
$ws.Range("L61").Value = 2066.6667
$ws.Range("M61").Value = -929.5
$ws.Range("N61").Value = -2490.6667

$ws.Range("H132").Value = 2582.0557
$ws.Range("I132").Value = 2650
$ws.Range("K132").Value = 7950
$ws.Range("M132").Value = -5420

$ws.Range("H136").Value = 1449.8889
$ws.Range("I136").Value = 1141.5
$ws.Range("J136").Value = 2066.6667
$ws.Range("K136").Value = 3424.5
$ws.Range("L136").Value = 6200.000100000001
$ws.Range("M136").Value = -874.5
$ws.Range("N136").Value = -11300.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

$ws.Range("H8").Value = 100
$ws.Range("I8").Value = 100
$ws.Range("K8").Value = 100
$ws.Range("M8").Value = 40

$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H19").Value = 498
$ws.Range("J19").Value = 498
$ws.Range("L19").Value = 498
$ws.Range("N19").Value = -844

$ws.Range("H82").Value = 31499.5
$ws.Range("I82").Value = 29999
$ws.Range("K82").Value = 29999
$ws.Range("M82").Value = -29616

$ws.Range("H85").Value = 31499.5
$ws.Range("I85").Value = 29999
$ws.Range("K85").Value = 29999
$ws.Range("M85").Value = -28673

$ws.Range("H134").Value = 8202.277
$ws.Range("I134").Value = 752.8
$ws.Range("J134").Value = 17514.125
$ws.Range("K134").Value = 2258.4
$ws.Range("L134").Value = 52542.375
$ws.Range("M134").Value = 276.6000000000004
$ws.Range("N134").Value = -57612.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 256.5
$ws.Range("J7").Value = 256.5
$ws.Range("L7").Value = 256.5
$ws.Range("N7").Value = -482.5

$ws.Range("H16").Value = 76924170
$ws.Range("I16").Value = 111112140
$ws.Range("J16").Value = 1225
$ws.Range("K16").Value = 111112140
$ws.Range("L16").Value = 1225
$ws.Range("M16").Value = -111111853
$ws.Range("N16").Value = -1799

$ws.Range("H22").Value = 364.2857
$ws.Range("I22").Value = 290
$ws.Range("J22").Value = 550
$ws.Range("K22").Value = 290
$ws.Range("L22").Value = 550
$ws.Range("M22").Value = 60
$ws.Range("N22").Value = -1250

$ws.Range("H31").Value = 1542.5834
$ws.Range("I31").Value = 1251.1
$ws.Range("K31").Value = 1251.1
$ws.Range("M31").Value = -956.0999999999999

$ws.Range("H34").Value = 1542.5834
$ws.Range("I34").Value = 1251.1
$ws.Range("K34").Value = 1251.1
$ws.Range("M34").Value = -1049.1

$ws.Range("H93").Value = 18600
$ws.Range("I93").Value = 4500
$ws.Range("K93").Value = 4500
$ws.Range("M93").Value = -2628

$ws.Range("H113").Value = 76924170
$ws.Range("I113").Value = 111112140
$ws.Range("J113").Value = 1225
$ws.Range("K113").Value = 111112140
$ws.Range("L113").Value = 1225
$ws.Range("M113").Value = -111109970
$ws.Range("N113").Value = -5565

$ws.Range("H132").Value = 14212.111
$ws.Range("I132").Value = 26728
$ws.Range("J132").Value = 4199.4
$ws.Range("K132").Value = 80184
$ws.Range("L132").Value = 12598.2
$ws.Range("M132").Value = -77654
$ws.Range("N132").Value = -17658.2

$ws.Range("H134").Value = 66670070
$ws.Range("I134").Value = 111115120
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 333345360
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -333342825
$ws.Range("N134").Value = -12570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 2700
$ws.Range("I125").Value = 2700
$ws.Range("K125").Value = 8100
$ws.Range("M125").Value = -3180

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 500.3846
$ws.Range("I107").Value = 659.0769
$ws.Range("J107").Value = 341.69232
$ws.Range("K107").Value = 659.0769
$ws.Range("L107").Value = 341.69232
$ws.Range("M107").Value = 1260.9231
$ws.Range("N107").Value = -4181.69232

$ws.Range("H113").Value = 1780
$ws.Range("I113").Value = 1117.4286
$ws.Range("J113").Value = 2707.6
$ws.Range("K113").Value = 1117.4286
$ws.Range("L113").Value = 2707.6
$ws.Range("M113").Value = 1052.5714
$ws.Range("N113").Value = -7047.6

$ws.Range("H132").Value = 2507.3076
$ws.Range("I132").Value = 2126.0527
$ws.Range("J132").Value = 3542.1428
$ws.Range("K132").Value = 6378.158100000001
$ws.Range("L132").Value = 10626.4284
$ws.Range("M132").Value = -3848.158100000001
$ws.Range("N132").Value = -15686.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2628.8572
$ws.Range("J46").Value = 2628.8572
$ws.Range("L46").Value = 2628.8572
$ws.Range("N46").Value = -3004.8572

$ws.Range("H51").Value = 10084
$ws.Range("J51").Value = 10084
$ws.Range("L51").Value = 10084
$ws.Range("N51").Value = -11040

$ws.Range("H61").Value = 1297
$ws.Range("J61").Value = 1824.75
$ws.Range("L61").Value = 1824.75
$ws.Range("N61").Value = -2228.75

$ws.Range("H75").Value = 30173
$ws.Range("J75").Value = 30173
$ws.Range("L75").Value = 30173
$ws.Range("N75").Value = -32045

$ws.Range("H78").Value = 30173
$ws.Range("J78").Value = 30173
$ws.Range("L78").Value = 90519
$ws.Range("N78").Value = -99879

$ws.Range("H113").Value = 1297
$ws.Range("J113").Value = 1824.75
$ws.Range("L113").Value = 1824.75
$ws.Range("N113").Value = -6164.75

$ws.Range("H132").Value = 24968.72
$ws.Range("I132").Value = 1460.7727
$ws.Range("J132").Value = 49596.094
$ws.Range("K132").Value = 4382.3181
$ws.Range("L132").Value = 148788.282
$ws.Range("M132").Value = -1852.3181
$ws.Range("N132").Value = -153848.282

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 8989.25
$ws.Range("J51").Value = 8989.25
$ws.Range("L51").Value = 8989.25
$ws.Range("N51").Value = -10009.25

$ws.Range("H132").Value = 3913.1333
$ws.Range("I132").Value = 3580.4
$ws.Range("J132").Value = 4578.6
$ws.Range("K132").Value = 10741.2
$ws.Range("L132").Value = 13735.8
$ws.Range("M132").Value = -8211.200000000001
$ws.Range("N132").Value = -18795.8
